# Update paths for test data.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("pages_with_simple_dates")
$ws2 = $wb.Worksheets.Item("pages_without_dates")

# --- sheet "pages_with_simple_dates": update the reconstruction-fact-sheet rows ---
# Row 7 (Spanish "Article" / "Revisión:" row): path + date label + ISO stamp change
# from the old "hoja-informativa-estres" page to the new reconstruction fact sheet page.
$ws1.Range("A7").Value = "espanol/hoja-informativa-reconstruccion"
$ws1.Range("E7").Value = "12 de marzo de 2016"
$ws1.Range("F7").Value = "2016-03-12T12:00:00Z"

# Row 6 (English "Article" / "Reviewed:" row): only the path changes.
$ws1.Range("A6").Value = "types/breast/reconstruction-fact-sheet"

# --- bold the header row on both sheets ---
$ws1.Range("A1:F1").Font.Bold = $true
$ws2.Range("A1:C1").Font.Bold = $true

# --- refresh the selection on both sheets to the header row ---
[void]$ws1.Range("A1:XFD1").Select()
[void]$ws2.Range("A1:XFD1").Select()

# Re-activate sheet 1 (the tab that is selected in the saved workbook).
$ws1.Activate()
